$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New FAQ content (Chinese + Malay) appended as rows 65-70 ----
$A65 = @"
什么是脐带血？
"@
$B65 = @"
脐带血是指在婴儿出生后以及脐带被剪断后仍残留在脐带中的血液。怀孕期间，脐带充当母亲和婴儿之间的生命线。在婴儿出生后，脐带中的脐带血可能为婴儿以及家庭成员提供希望。脐带血含有一种被称为造血干细胞（Haematopoietic Stem Cells，HSCs）的具有生命拯救作用的干细胞。
"@
$A66 = @"
什么是脐带血干细胞？
"@
$B66 = @"
“脐带血干细胞也被称为造血干细胞（Haematopoietic Stem Cells，HSCs），它们负责补充血液并再生免疫系统。这些干细胞可在婴儿出生后的脐带中找到，并具有分化为血液中各种细胞类型的独特能力，如下图所示： 红细胞 - 运输氧气 白细胞 - 产生抗体并抵抗细菌 血小板 - 协助血液凝固 ”
"@
$A67 = @"
为什么要保存宝宝的脐带血干细胞？
"@
$B67 = @"
“保存宝宝的脐带血干细胞有几个优点，例如：
确保自体移植的匹配（供体和受体是同一人）
拥有随时可用的储存的造血干细胞（HSCs），而不是进行国家或国际搜索，这在紧急情况下可能既昂贵又耗时
对于自体移植，移植组织攻击患者自身组织的移植物反应（Graft vs. Host Disease，GvHD）的风险较低
无创收集程序对母亲和婴儿都是无痛且无风险的
与其他干细胞来源（例如骨髓）相比，脐带血干细胞更年轻，移植率更高，对组织不匹配更容忍”
"
"@
$A68 = @"
Apakah darah tali pusat? 
"@
$B68 = @"
Darah tali pusat ialah darah yang kekal dalam tali pusat selepas kelahiran bayi dan selepas tali pusat dipotong. Semasa kehamilan, tali pusat berfungsi sebagai talian hayat antara ibu dan anak. Selepas kelahiran bayi, darah tali pusat yang terdapat dalam tali pusat boleh memberikan harapan untuk bayi dan ahli keluarga. Darah tali pusat mengandungi sumber yang kaya dengan sel punca penyelamat yang dipanggil Sel Punca Hematopoietik (HSC).
"@
$A69 = @"
Apakah sel punca darah tali pusat?
"@
$B69 = @"
"Sel punca darah tali pusat juga dikenali sebagai Sel Punca Hematopoietik (HSC), yang bertanggungjawab untuk menyokong semula darah dan memulihkan sistem imun. Sel punca ini boleh dijumpai dalam tali pusat selepas kelahiran bayi dan mempunyai keupayaan unik untuk berbeza menjadi pelbagai jenis sel yang terdapat dalam darah seperti yang digambarkan dalam gambarajah di bawah ini:
Sel darah merah - yang mengangkut oksigen
Sel darah putih - yang menghasilkan antibodi dan memerangi bakteria
Platelet - yang membantu dalam pembekuan darah
"
"@
$A70 = @"
Mengapa saya perlu menyimpan sel punca darah tali pusat bayi saya?
"@
$B70 = @"
"Terdapat beberapa kelebihan menyimpan sel punca darah tali pusat bayi anda, seperti:
Padanan dijamin untuk transplantasi autologous (di mana donor dan penerima adalah individu yang sama)
Mempunyai bekalan Sel Punca Hematopoietik (HSC) yang disimpan yang sedia ada berbanding menjalankan carian kebangsaan atau antarabangsa, yang boleh mahal dan memakan masa semasa situasi yang memerlukan tindakan segera
Risiko yang rendah untuk Penyakit Graft vs. Host (GvHD) bagi transplantasi autologous, di mana tisu yang dipindahkan menyerang tisu sendiri pesakit
Prosedur pengumpulan yang tidak invasif yang tidak menyakitkan dan bebas risiko kepada ibu dan bayi
Sel punca darah tali pusat lebih muda, mempunyai kadar pengapitan yang lebih tinggi, dan lebih toleran terhadap tidak sepadan tisu berbanding sumber sel punca lain, seperti sumsum tulang"
"
"@

# ---- Write values into the worksheet cells ----
$ws.Range("A65").Value = $A65
$ws.Range("B65").Value = $B65
$ws.Range("A66").Value = $A66
$ws.Range("B66").Value = $B66
$ws.Range("A67").Value = $A67
$ws.Range("B67").Value = $B67
$ws.Range("A68").Value = $A68
$ws.Range("B68").Value = $B68
$ws.Range("A69").Value = $A69
$ws.Range("B69").Value = $B69
$ws.Range("A70").Value = $A70
$ws.Range("B70").Value = $B70

# ---- Apply new font style (Segoe UI, 8pt, color #374151) ----
# Build the style from scratch on A65, then propagate via copy/paste-format
# (keeps the style table minimal / matches a single new cellXfs entry)
$styleSource = $ws.Range("A65")
$styleSource.Style = "Normal"
$styleSource.Font.Name = "Segoe UI"
$styleSource.Font.Size = 8
$styleSource.Font.Color = 5325111

$styleSource.Copy()
foreach ($addr in @("A66","A67","A68","A69","A70","B65","B66","B68")) {
    $ws.Range($addr).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# ---- Row heights for the long-answer rows ----
$ws.Rows.Item(67).RowHeight = 116
$ws.Rows.Item(69).RowHeight = 116
$ws.Rows.Item(70).RowHeight = 159.5

# ---- Update the active view to match the final cursor position ----
$excel.ActiveWindow.ScrollRow = 57
$ws.Range("B71").Select() | Out-Null
